$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feature Importance")

# New ordering (row index 2..9) with updated p-value (B) and objective value (D).
# Importance (C) values are unchanged per-label.
$data = @(
    @{ Row = 2; Name = "Tx";     B = 0.2893169989624623;  C = 0.3943872552251285;  D = 1 },
    @{ Row = 3; Name = "Tn";     B = 0.1813592958955645;  C = 0.1816243256743412;  D = 0.4237924090919176 },
    @{ Row = 4; Name = "RH_avg"; B = 0.2506164845172436;  C = 0.1061353905077548;  D = 0.2193522015595606 },
    @{ Row = 5; Name = "Tavg";   B = 0.06065280242151002; C = 0.1052214880401069;  D = 0.2168771579026596 },
    @{ Row = 6; Name = "ss";     B = 0.09340294057062298; C = 0.07972246640130111; D = 0.1478203449511246 },
    @{ Row = 7; Name = "ff_avg"; B = 0.07646911220056765; C = 0.06001532725417928; D = 0.09444919051351454 },
    @{ Row = 8; Name = "ff_x";   B = 0.07142938879604899; C = 0.0477535007499524;  D = 0.06124153763408323 },
    @{ Row = 9; Name = "RR";     B = 0.1304048800547975;  C = 0.02514024614723597; D = 0 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Name
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
}
